# bug fixed and some new functions
# Update a handful of stat values on the (only) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 431016551261405195)
$ws.Range("B2").Value = 3394548834   # was 3394326834
$ws.Range("G2").Value = 0            # was 11

# Row 3 (A3 = 568754491101544448)
$ws.Range("D3").Value = 360          # was 350
$ws.Range("I3").Value = 2            # was 1

# Row 4 (A4 = 539455628159090688)
$ws.Range("B4").Value = -2000        # was 0
